$d = $word.ActiveDocument

# --- Simple, unambiguous whole-document text replacements ---
$d.Content.Find.Execute("numline", $true, $false, $false, $false, $false, $true, 1, $false, "ความรู้ทั่วไป", 2)
$d.Content.Find.Execute("xx999", $true, $false, $false, $false, $false, $true, 1, $false, "GN002", 2)
$d.Content.Find.Execute("มิถุนายน", $true, $false, $false, $false, $false, $true, 1, $false, "พฤษภาคม", 2)
$d.Content.Find.Execute("01:00", $true, $false, $false, $false, $false, $true, 1, $false, "10:00", 2)
$d.Content.Find.Execute("11:00", $true, $false, $false, $false, $false, $true, 1, $false, "12:00", 2)

# --- Cell-scoped replacements for ambiguous "1" values in the first (header) table ---
# (Cell addressing is by grid column, so the columns below fall inside the
#  correct merged cell for each field.)
$examTable = $d.Tables(1)

# "14" -> "21" (exam day, row 4)
$examTable.Cell(4, 2).Range.Find.Execute("14", $true, $false, $false, $false, $false, $true, 1, $false, "21", 2)

# "1" -> "3" (credit units, row 1)
$examTable.Cell(1, 22).Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2)

# "1" -> "2" (semester number, row 5)
$examTable.Cell(5, 6).Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2)

# --- Remove the "คำสั่ง" instructions row from the second table ---
$instructionsTable = $d.Tables(2)
$instructionsTable.Rows(1).Delete()
